$d = $word.ActiveDocument

$replacements = @(
    @("348×6=2088", "863×3=2589"),
    @("966×2=1932", "552×8=4416"),
    @("456×9=4104", "734×5=3670"),
    @("209×8=1672", "491×5=2455"),
    @("333×6=1998", "574×6=3444"),
    @("559×9=5031", "344×4=1376"),
    @("327×3=981",  "979×7=6853"),
    @("511×8=4088", "221×7=1547"),
    @("376×9=3384", "766×2=1532"),
    @("871×9=7839", "608×4=2432"),
    @("272×6=1632", "115×5=575"),
    @("983×4=3932", "611×6=3666"),
    @("843×5=4215", "897×6=5382"),
    @("497×7=3479", "567×3=1701"),
    @("689×6=4134", "323×5=1615"),
    @("962×8=7696", "138×9=1242"),
    @("489×3=1467", "924×6=5544"),
    @("287×8=2296", "167×5=835"),
    @("647×6=3882", "267×4=1068"),
    @("885×8=7080", "774×2=1548"),
    @("837×4=3348", "853×4=3412"),
    @("413×9=3717", "159×4=636"),
    @("926×9=8334", "873×6=5238"),
    @("178×6=1068", "345×6=2070"),
    @("528×4=2112", "427×4=1708")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
